$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" "29.833.22"
Set-TextValue "E2" "  +2.47%  "
Set-TextValue "D3" "1.857.99"
Set-TextValue "E3" "  +1.81%  "
Set-TextValue "D4" "0.9993"
Set-TextValue "E4" "  +0.03%  "
Set-TextValue "D5" "244.23"
Set-TextValue "E5" "  +1.18%  "
Set-TextValue "D6" "0.6423"
Set-TextValue "E6" "  +3.80%  "
Set-TextValue "D7" "0.9994"
Set-TextValue "E7" "  -0.01%  "
Set-TextValue "B8" "Dogecoin"
Set-TextValue "C8" "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue "D8" "0.07505"
Set-TextValue "E8" "  +2.42%  "
Set-TextValue "B9" "Cardano"
Set-TextValue "C9" "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue "D9" "0.2967"
Set-TextValue "E9" "  +2.10%  "
Set-TextValue "B10" "Solana"
Set-TextValue "C10" "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D10" "24.39"
Set-TextValue "E10" "  +5.95%  "
Set-TextValue "B11" "TRON"
Set-TextValue "C11" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D11" "0.07662"
Set-TextValue "E11" "  -0.09%  "
Set-TextValue "B12" "WrappedEther"
Set-TextValue "C12" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D12" "1.864.59"
Set-TextValue "E12" "  +2.15%  "
Set-TextValue "B13" "Polkadot"
Set-TextValue "C13" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D13" "5.063"
Set-TextValue "E13" "  +2.55%  "
Set-TextValue "B14" "Polygon"
Set-TextValue "C14" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D14" "0.6907"
Set-TextValue "E14" "  +4.29%  "
Set-TextValue "B15" "Litecoin"
Set-TextValue "C15" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D15" "84.22"
Set-TextValue "E15" "  +2.83%  "
Set-TextValue "B16" "ShibaInu"
Set-TextValue "C16" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D16" "0.000009583"
Set-TextValue "E16" "  +7.43%  "
Set-TextValue "B17" "Uniswap"
Set-TextValue "C17" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D17" "6.071"
Set-TextValue "E17" "  +4.20%  "
Set-TextValue "B18" "WrappedBTC"
Set-TextValue "C18" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue "D18" "29.821.62"
Set-TextValue "E18" "  +2.52%  "
Set-TextValue "B19" "WrappedliquidstakedEther2.0"
Set-TextValue "C19" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D19" "2.111.40"
Set-TextValue "E19" "  +1.83%  "
Set-TextValue "B20" "BitcoinCash"
Set-TextValue "C20" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D20" "238.48"
Set-TextValue "E20" "  -0.19%  "
Set-TextValue "B21" "Avalanche"
Set-TextValue "C21" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D21" "12.68"
Set-TextValue "E21" "  +2.31%  "
Set-TextValue "B22" "Dai"
Set-TextValue "C22" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D22" "0.9998"
Set-TextValue "E22" "  +0.06%  "
Set-TextValue "B23" "Chainlink"
Set-TextValue "C23" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D23" "7.449"
Set-TextValue "E23" "  +1.63%  "
Set-TextValue "B24" "BinanceUSD"
Set-TextValue "C24" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D24" "1.000"
Set-TextValue "E24" "  +0.00%  "
Set-TextValue "B25" "Monero"
Set-TextValue "C25" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D25" "158.65"
Set-TextValue "E25" "  +0.60%  "
Set-TextValue "B26" "Stellar"
Set-TextValue "C26" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D26" "0.1430"
Set-TextValue "E26" "  +0.76%  "
Set-TextValue "B27" "Cosmos"
Set-TextValue "C27" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D27" "8.534"
Set-TextValue "E27" "  +0.81%  "
Set-TextValue "B28" "EthereumClassic"
Set-TextValue "C28" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D28" "17.95"
Set-TextValue "E28" "  +1.82%  "
Set-TextValue "B29" "Hedera"
Set-TextValue "C29" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D29" "0.06140"
Set-TextValue "E29" "  +3.53%  "
Set-TextValue "B30" "PancakeSwap"
Set-TextValue "C30" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D30" "1.495"
Set-TextValue "E30" "  +0.98%  "
Set-TextValue "B31" "Toncoin"
Set-TextValue "C31" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D31" "1.271"
Set-TextValue "E31" "  +5.88%  "
Set-TextValue "B32" "Filecoin"
Set-TextValue "C32" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D32" "4.147"
Set-TextValue "E32" "  +1.81%  "
Set-TextValue "B33" "InternetComputer(DFINITY)"
Set-TextValue "C33" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D33" "4.105"
Set-TextValue "E33" "  +1.20%  "
Set-TextValue "B34" "LidoDAOToken"
Set-TextValue "C34" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D34" "1.878"
Set-TextValue "E34" "  +0.80%  "
Set-TextValue "B35" "ARBITRUM"
Set-TextValue "C35" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D35" "1.171"
Set-TextValue "E35" "  +3.45%  "
Set-TextValue "B36" "ImmutableX"
Set-TextValue "C36" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D36" "0.7325"
Set-TextValue "E36" "  +0.30%  "
Set-TextValue "B37" "HuobiToken"
Set-TextValue "C37" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D37" "2.614"
Set-TextValue "E37" "  +0.53%  "
Set-TextValue "B38" "MXToken"
Set-TextValue "C38" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D38" "2.853"
Set-TextValue "E38" "  +0.30%  "
Set-TextValue "B39" "VeChain"
Set-TextValue "C39" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D39" "0.01794"
Set-TextValue "E39" "  +2.49%  "
Set-TextValue "B40" "Maker"
Set-TextValue "C40" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D40" "1.214.91"
Set-TextValue "E40" "  +0.38%  "
Set-TextValue "B41" "TrustWalletToken"
Set-TextValue "C41" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D41" "0.9273"
Set-TextValue "E41" "  +0.43%  "
Set-TextValue "B42" "FraxShare"
Set-TextValue "C42" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D42" "6.186"
Set-TextValue "E42" "  -1.02%  "
Set-TextValue "B43" "PaxDollar"
Set-TextValue "C43" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D43" "1.000"
Set-TextValue "E43" "  +0.04%  "
Set-TextValue "B44" "RocketPoolETH"
Set-TextValue "C44" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D44" "2.022.77"
Set-TextValue "E44" "  +2.30%  "
Set-TextValue "B45" "Quant"
Set-TextValue "C45" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D45" "102.39"
Set-TextValue "E45" "  +0.73%  "
Set-TextValue "B46" "Aave"
Set-TextValue "C46" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D46" "66.30"
Set-TextValue "E46" "  +2.44%  "
Set-TextValue "B47" "BabyDogeCoin"
Set-TextValue "C47" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D47" "0.00000000121"
Set-TextValue "E47" "  -2.27%  "
Set-TextValue "B48" "EnergySwap"
Set-TextValue "C48" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D48" "9.246"
Set-TextValue "E48" "  +1.81%  "
Set-TextValue "D49" "0.4077"
Set-TextValue "E49" "  +1.78%  "
Set-TextValue "B50" "Cronos"
Set-TextValue "C50" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D50" "0.05820"
Set-TextValue "E50" "  +1.17%  "
Set-TextValue "B51" "RenderToken"
Set-TextValue "C51" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D51" "1.658"
Set-TextValue "E51" "  +5.12%  "
